$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 10; this shifts existing rows 10..50 down to 11..51
$ws.Rows.Item(10).Insert()

# Populate the newly inserted row 10 with the new weekly data record
$ws.Cells.Item(10, 1).Value = 1
$ws.Cells.Item(10, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(10, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(10, 4).Value = 44592
$ws.Cells.Item(10, 5).Value = 15
$ws.Cells.Item(10, 6).Value = 100112009
$ws.Cells.Item(10, 7).Value = "Acelga"
$ws.Cells.Item(10, 8).Value = "Sin especificar"
$ws.Cells.Item(10, 9).Value = "Primera"
$ws.Cells.Item(10, 10).Value = 200
$ws.Cells.Item(10, 11).Value = 1800
$ws.Cells.Item(10, 12).Value = 2000
$ws.Cells.Item(10, 13).Value = 1900
$ws.Cells.Item(10, 14).Value = "$/atado 2,5 a 3 kilos"
$ws.Cells.Item(10, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(10, 16).Value = 633
$ws.Cells.Item(10, 17).Value = 3
$ws.Cells.Item(10, 18).Value = "Hortaliza"
